$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Goals")

# --- Row 30 (Id 29): "Diving" -> "Save the Apple Reward", Difficulty 1 -> 14, clear Subtype ---
$ws.Range("B30").Value = "Save the Apple Reward"
$ws.Range("B30").Font.Name = "Arial"
$ws.Range("B30").Font.Size = 11

$ws.Range("C30").Value = 14
$ws.Range("C30").Font.Name = "Arial"
$ws.Range("C30").Font.Size = 11
$ws.Range("C30").HorizontalAlignment = -4152

$ws.Range("E30").Value = ""

# --- Row 31 (Id 30): "Vine Swinging" -> "Castle Lanky Tower", Difficulty 1 -> 14, clear Subtype ---
$ws.Range("B31").Value = "Castle Lanky Tower"
$ws.Range("B31").Font.Name = "Arial"
$ws.Range("B31").Font.Size = 11

$ws.Range("C31").Value = 14
$ws.Range("C31").Font.Name = "Arial"
$ws.Range("C31").Font.Size = 11
$ws.Range("C31").HorizontalAlignment = -4152

$ws.Range("E31").Value = ""

# --- Row 44 (Id 43): "Disable the Blastomatic" -> "1 Company Coin", Difficulty 20 -> 12, clear Type ---
$ws.Range("B44").Value = "1 Company Coin"
$ws.Range("B44").Font.Name = "Arial"
$ws.Range("B44").Font.Size = 11

$ws.Range("C44").Value = 12
$ws.Range("C44").Font.Name = "Arial"
$ws.Range("C44").Font.Size = 11
$ws.Range("C44").HorizontalAlignment = -4152

$ws.Range("D44").Value = ""
